$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Inserts a new paragraph with text $newText immediately after the (unique)
# paragraph whose text equals $searchText. Pass "" for $newText to insert a
# blank paragraph. The new paragraph never inherits paragraph-level styling
# (e.g. heading styles) because it is built via InsertBefore on a collapsed
# range positioned just after the anchor paragraph's mark.
function Insert-ParaAfter($searchText, $newText) {
    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $searchText) {
            $r = $p.Range
            $r.Collapse(0)
            $r.InsertBefore($newText + [char]13)
            return
        }
    }
}

# 1. Update the "Created" timestamp.
Replace-Text "Created April 26, 2025 at 10:46:01" "Created April 27, 2025 at 15:53:26"

# 2./3. Right after "System Overview": a blank paragraph, the new descriptive
#    sentence, another blank paragraph, then "* Date: 04-06-2025" (the
#    original "Date: 04-06-2025" paragraph, edited in place).
#    Insert the sentence right after the heading first, then a blank before
#    it (also anchored on the heading), and finally a blank after it
#    (anchored on the now-unique sentence text).
Insert-ParaAfter "System Overview" "The system being patched is a Debian-based operating system named `"kb322-18`". The system was last updated on April 6, 2025, at 4:02 PM, and its current status is as follows:"
Insert-ParaAfter "System Overview" ""
Insert-ParaAfter "The system being patched is a Debian-based operating system named `"kb322-18`". The system was last updated on April 6, 2025, at 4:02 PM, and its current status is as follows:" ""
Replace-Text "Date: 04-06-2025" "* Date: 04-06-2025"

# 4. Prefix the remaining System Overview detail lines with "* ".
Replace-Text "Time: 16:02:09" "* Time: 16:02:09"
Replace-Text "OS Name: kb322-18" "* OS Name: kb322-18"
Replace-Text "OS Version: #1 SMP PREEMPT_DYNAMIC Debian 6.1.129-1 (2025-03-06)" "* OS Version: #1 SMP PREEMPT_DYNAMIC Debian 6.1.129-1 (2025-03-06)"
Replace-Text "Computer Name: kb322-18" "* Computer Name: kb322-18"
Replace-Text "IP Address: 140.160.138.147" "* IP Address: 140.160.138.147"

# 5. Blank paragraph after "Patch Status Summary" heading.
Insert-ParaAfter "Patch Status Summary" ""

# 6. Blank paragraph after "Compliance with RMF Controls" heading.
Insert-ParaAfter "Compliance with RMF Controls" ""

# 7. Replace the RMF controls paragraph text.
Replace-Text "To ensure compliance, it is essential to have a process in place for identifying and reporting vulnerabilities. The current system does not have any identified patches that require remediation." "There are no patches to apply, therefore, there is no need for immediate corrective action."

# 8. Blank paragraph after "Recommended next steps" heading.
Insert-ParaAfter "Recommended next steps" ""

# 9. Replace the three recommendation lines.
Replace-Text "Review the current patch status and schedule any necessary patch deployments." "* Review the system for any upcoming security patches."
Replace-Text "Update documentation should be reviewed and updated to reflect the current patch status." "* Schedule patch deployments as needed."
Replace-Text "Schedule regular review and assessment of updates to ensure the system remains secure." "* Maintain accurate records of all update documentation."

# 10. Two blank paragraphs after "Risk Assessment" heading.
Insert-ParaAfter "Risk Assessment" ""
Insert-ParaAfter "Risk Assessment" ""

# 11. Replace the final risk-assessment paragraph text.
Replace-Text "Since there are no pending updates, the risk level is considered low. There is no potential impact on the system's security at this time. However, it is still essential to regularly review and assess updates to ensure the system remains secure." "There are currently no pending updates available. Therefore, there is no potential risk to report at this time. However, it is recommended to regularly review and assess system updates to ensure the system remains secure."
